$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.340.67'
$ws.Range("E2").Value = '  +0.46%  '
$ws.Range("D3").Value = '1.932.26'
$ws.Range("E3").Value = '  +0.40%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.22%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '250.67'
$ws.Range("E5").Value = '  +1.98%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.7166'
$ws.Range("E6").Value = '  -0.29%  '
$ws.Range("E7").Value = '  +0.26%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3263'
$ws.Range("E8").Value = '  +0.66%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '27.51'
$ws.Range("E9").Value = '  +3.96%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07174'
$ws.Range("E10").Value = '  +4.83%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.7989'
$ws.Range("E11").Value = '  +0.28%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08083'
$ws.Range("E12").Value = '  +2.04%  '
$ws.Range("D13").Value = '1.928.59'
$ws.Range("E13").Value = '  +0.21%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.413'
$ws.Range("E14").Value = '  +0.47%  '
$ws.Range("E15").Value = '  +0.20%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.76'
$ws.Range("E16").Value = '  +1.93%  '
$ws.Range("D17").Value = '30.319.66'
$ws.Range("E17").Value = '  +0.38%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '251.32'
$ws.Range("E18").Value = '  -3.40%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000008117'
$ws.Range("E19").Value = '  +2.18%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.789'
$ws.Range("E20").Value = '  -0.36%  '
$ws.Range("D21").Value = '2.184.31'
$ws.Range("E21").Value = '  +0.45%  '
$ws.Range("E22").Value = '  +0.27%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.005'
$ws.Range("E23").Value = '  +0.64%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.910'
$ws.Range("E24").Value = '  +0.82%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.709'
$ws.Range("E25").Value = '  +0.61%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '165.77'
$ws.Range("E26").Value = '  +3.47%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.20'
$ws.Range("E27").Value = '  +1.64%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.324'
$ws.Range("E28").Value = '  +2.97%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1284'
$ws.Range("E29").Value = '  -3.96%  '
$ws.Range("E30").Value = '  +0.48%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.543'
$ws.Range("E31").Value = '  +0.14%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.425'
$ws.Range("E32").Value = '  +0.28%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.193'
$ws.Range("E33").Value = '  +0.25%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05202'
$ws.Range("E34").Value = '  +3.04%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.263'
$ws.Range("E35").Value = '  +6.25%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7461'
$ws.Range("E36").Value = '  +1.27%  '
$ws.Range("E37").Value = '  +1.26%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01958'
$ws.Range("E38").Value = '  +1.16%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.800'
$ws.Range("E39").Value = '  -0.25%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '78.92'
$ws.Range("E40").Value = '  -1.41%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.442'
$ws.Range("E41").Value = '  -1.19%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4519'
$ws.Range("E42").Value = '  +1.73%  '
$ws.Range("E43").Value = '  +0.78%  '
$ws.Range("E44").Value = '  +0.19%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8407'
$ws.Range("E45").Value = '  +1.23%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '101.80'
$ws.Range("E46").Value = '  -0.76%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.819'
$ws.Range("E47").Value = '  +0.74%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.399'
$ws.Range("E48").Value = '  +1.85%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '36.59'
$ws.Range("E49").Value = '  +1.12%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06072'
$ws.Range("E50").Value = '  +2.77%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4171'
$ws.Range("E51").Value = '  +1.78%  '
